# Auto-generated edit script: refresh cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.126.49"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "3.402.41"
$ws.Range("E3").Value = "  -1.41%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "573.65"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.83%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "142.52"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("D7").Value = "3.402.98"
$ws.Range("E7").Value = "  -1.44%  "
$ws.Range("E9").Value = "  -0.69%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "7.61"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("E11").Value = "  -2.30%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.398"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("D13").Value = "3.983.21"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("E14").Value = "  +2.10%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "28.09"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").Value = "3.403.35"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("D18").Value = "61.154.16"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("E19").Value = "  -3.49%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.87"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.28%  "
$ws.Range("E21").Value = "  -4.85%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "384.03"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -4.90%  "
$ws.Range("E23").Value = "  -1.50%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "74.52"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").Value = "  -5.24%  "
$ws.Range("D27").Value = "3.538.60"
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("E28").Value = "  -1.66%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("E30").Value = "  -3.04%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "8.02"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.68%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.16"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.30%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.42"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.91%  "
$ws.Range("E34").Value = "  -0.08%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "23.52"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.97%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "7.02"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.59%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "167.79"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").Value = "3.432.75"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("E39").Value = "  -2.76%  "
$ws.Range("E40").Value = "  -5.39%  "
$ws.Range("E41").Value = "  -2.14%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "27.25"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("E44").Value = "  -0.03%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "4.44"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.81%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.67"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.54%  "
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("D48").Value = "2.487.36"
$ws.Range("E48").Value = "  -4.85%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "6.83"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.97%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "22.95"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("E51").Value = "  +1.13%  "
